{"js": "// Update the date heading (first paragraph in the body, before the table).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items[0].insertText(\"2025-07-13 Sunday\", \"Replace\");\n\n// Update every answer cell in the practice-problems table with the new\n// expressions, preserving row-major (row 1..20, col 1..5) order.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.values = [\n  [\n    \"7+24=\",\n    \"96-9=\",\n    \"95-37=\",\n    \"91-63=\",\n    \"94-88=\"\n  ],\n  [\n    \"25+48=\",\n    \"68-49=\",\n    \"80-17=\",\n    \"9+13=\",\n    \"17+55=\"\n  ],\n  [\n    \"33-18=\",\n    \"25+57=\",\n    \"72-49=\",\n    \"12-7=\",\n    \"15+19=\"\n  ],\n  [\n    \"93-79=\",\n    \"90-24=\",\n    \"60-8=\",\n    \"94-59=\",\n    \"85-76=\"\n  ],\n  [\n    \"91-9=\",\n    \"62+29=\",\n    \"86-47=\",\n    \"92-79=\",\n    \"19+8=\"\n  ],\n  [\n    \"39+44=\",\n    \"87+9=\",\n    \"66+9=\",\n    \"37+44=\",\n    \"38+9=\"\n  ],\n  [\n    \"54-15=\",\n    \"9+13=\",\n    \"40-12=\",\n    \"6+78=\",\n    \"44+9=\"\n  ],\n  [\n    \"84-29=\",\n    \"81-25=\",\n    \"90-24=\",\n    \"60-35=\",\n    \"29+45=\"\n  ],\n  [\n    \"38+24=\",\n    \"4+49=\",\n    \"22+49=\",\n    \"90-18=\",\n    \"96-37=\"\n  ],\n  [\n    \"69+25=\",\n    \"17+9=\",\n    \"28+39=\",\n    \"72-65=\",\n    \"93-26=\"\n  ],\n  [\n    \"4+79=\",\n    \"33-29=\",\n    \"67-38=\",\n    \"74-39=\",\n    \"90-8=\"\n  ],\n  [\n    \"69+13=\",\n    \"72-29=\",\n    \"94-25=\",\n    \"87-78=\",\n    \"18+43=\"\n  ],\n  [\n    \"84-75=\",\n    \"73-29=\",\n    \"51-14=\",\n    \"86-49=\",\n    \"61-19=\"\n  ],\n  [\n    \"35+49=\",\n    \"75-26=\",\n    \"77+18=\",\n    \"35-19=\",\n    \"70-36=\"\n  ],\n  [\n    \"19+27=\",\n    \"77-68=\",\n    \"51-28=\",\n    \"73-8=\",\n    \"7+89=\"\n  ],\n  [\n    \"50-47=\",\n    \"79+12=\",\n    \"27+69=\",\n    \"41-29=\",\n    \"8+77=\"\n  ],\n  [\n    \"48+24=\",\n    \"18+26=\",\n    \"16+9=\",\n    \"82-6=\",\n    \"47+16=\"\n  ],\n  [\n    \"13+78=\",\n    \"61-14=\",\n    \"19+68=\",\n    \"37+18=\",\n    \"23+9=\"\n  ],\n  [\n    \"6+77=\",\n    \"92-16=\",\n    \"16+15=\",\n    \"51-7=\",\n    \"71-53=\"\n  ],\n  [\n    \"52+19=\",\n    \"65-58=\",\n    \"25+46=\",\n    \"43-7=\",\n    \"56-48=\"\n  ]\n];\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date heading paragraph\n$d.Paragraphs.Item(1).Range.Text = \"2025-07-13 Sunday\"\n\n# Update each answer cell in the practice table (row-major order)\n$answers = @(\n    \"7+24=\",\n    \"96-9=\",\n    \"95-37=\",\n    \"91-63=\",\n    \"94-88=\",\n    \"25+48=\",\n    \"68-49=\",\n    \"80-17=\",\n    \"9+13=\",\n    \"17+55=\",\n    \"33-18=\",\n    \"25+57=\",\n    \"72-49=\",\n    \"12-7=\",\n    \"15+19=\",\n    \"93-79=\",\n    \"90-24=\",\n    \"60-8=\",\n    \"94-59=\",\n    \"85-76=\",\n    \"91-9=\",\n    \"62+29=\",\n    \"86-47=\",\n    \"92-79=\",\n    \"19+8=\",\n    \"39+44=\",\n    \"87+9=\",\n    \"66+9=\",\n    \"37+44=\",\n    \"38+9=\",\n    \"54-15=\",\n    \"9+13=\",\n    \"40-12=\",\n    \"6+78=\",\n    \"44+9=\",\n    \"84-29=\",\n    \"81-25=\",\n    \"90-24=\",\n    \"60-35=\",\n    \"29+45=\",\n    \"38+24=\",\n    \"4+49=\",\n    \"22+49=\",\n    \"90-18=\",\n    \"96-37=\",\n    \"69+25=\",\n    \"17+9=\",\n    \"28+39=\",\n    \"72-65=\",\n    \"93-26=\",\n    \"4+79=\",\n    \"33-29=\",\n    \"67-38=\",\n    \"74-39=\",\n    \"90-8=\",\n    \"69+13=\",\n    \"72-29=\",\n    \"94-25=\",\n    \"87-78=\",\n    \"18+43=\",\n    \"84-75=\",\n    \"73-29=\",\n    \"51-14=\",\n    \"86-49=\",\n    \"61-19=\",\n    \"35+49=\",\n    \"75-26=\",\n    \"77+18=\",\n    \"35-19=\",\n    \"70-36=\",\n    \"19+27=\",\n    \"77-68=\",\n    \"51-28=\",\n    \"73-8=\",\n    \"7+89=\",\n    \"50-47=\",\n    \"79+12=\",\n    \"27+69=\",\n    \"41-29=\",\n    \"8+77=\",\n    \"48+24=\",\n    \"18+26=\",\n    \"16+9=\",\n    \"82-6=\",\n    \"47+16=\",\n    \"13+78=\",\n    \"61-14=\",\n    \"19+68=\",\n    \"37+18=\",\n    \"23+9=\",\n    \"6+77=\",\n    \"92-16=\",\n    \"16+15=\",\n    \"51-7=\",\n    \"71-53=\",\n    \"52+19=\",\n    \"65-58=\",\n    \"25+46=\",\n    \"43-7=\",\n    \"56-48=\"\n)\n\n$t = $d.Tables.Item(1)\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\n$i = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $t.Cell($r, $c).Range.Text = $answers[$i]\n        $i++\n    }\n}\n"}
